# DISCOVERYACCESS-4992: Add Fine Arts > Artists' Books collection to location facet mapping
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 16 (Excel shifts existing row 16.. down to 17..)
$ws.Rows.Item(16).Insert()

# Populate the new row with the Artist sub-collection mapping for Fine Arts Library
$ws.Cells.Item(16, 1).Value = "Fine Arts Library (B56 Sibley Hall)"
$ws.Cells.Item(16, 3).Value = "Artist"
$ws.Cells.Item(16, 5).Value = "Fine Arts Library > Artists' Books"

# Match the saved selection state
$ws.Range("E18").Select() | Out-Null
